$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# UniqueJobs counts refreshed for the Developer-keyword rows (stored as
# text, same as the rest of the UniqueJobs column, so prefix with an
# apostrophe to stop Excel from auto-converting to a number).
$ws.Range("G2").Value = "'696"
$ws.Range("G4").Value = "'1,955"

# Row 6 (previously "Entry Level Marketing" / DFW Brands) now holds the
# "Facilities Engineer" / National Security Agency posting data, fixing
# the duplicate UniqueJobs entry that used to live on row 7.
$ws.Range("B6").Value = "Facilities Engineer - Electrical - Entry/Experienced Level (NSAW and NSAH)"
$ws.Range("C6").Value = "National Security Agency (NSA)"
$ws.Range("D6").Value = "Fort Meade, MD"
$ws.Range("G6").Value = "'9"
$ws.Range("H6").Value = "12 days ago"
$ws.Range("L6").Value = "https://www.careerbuilder.com/job/J3Q4ML75ZRJFH5DYB94"

# Row 7 now holds the "Entry Level Marketing" / DFW Brands posting data
# (re-scraped, PostedDate bumped from "1 day ago" to "2 days ago"), and no
# longer duplicates the UniqueJobs value.
$ws.Range("B7").Value = "Entry Level Marketing"
$ws.Range("C7").Value = "DFW Brands"
$ws.Range("D7").Value = "Dallas, TX"
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = "2 days ago"
$ws.Range("L7").Value = "https://www.careerbuilder.com/job/J2W5L86NHXHH00GX3ZS"
